# Apply crypto price/volume updates (and a TrustWalletToken/TheSandbox row swap)
# as captured in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.325.01'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.875.38'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''235.09'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").Value = '''0.9998'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '''0.4697'
$ws.Range("E7").Value = '  +0.53%  '
$ws.Range("D8").Value = '''0.2874'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").Value = '''21.74'
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '''0.07953'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '''96.74'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = '1.877.60'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '''0.6952'
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Value = '''270.83'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '30.348.24'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '''14.04'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("D19").Value = '''0.000007749'
$ws.Range("E19").Value = '  +5.08%  '
$ws.Range("D20").Value = '''0.9997'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '2.122.94'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = '''0.9995'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '''5.272'
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("D24").Value = '''6.216'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = '''9.401'
$ws.Range("E25").Value = '  +2.09%  '
$ws.Range("D26").Value = '''167.50'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").Value = '''18.93'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").Value = '''1.953'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").Value = '''1.368'
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").Value = '''0.09891'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '''4.345'
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").Value = '''4.069'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '''0.04734'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = '''1.136'
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").Value = '''0.7030'
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").Value = '''0.01875'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '''2.810'
$ws.Range("E39").Value = '  +6.88%  '
$ws.Range("D40").Value = '''6.206'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").Value = '''72.20'
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").Value = '''1.960'
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.8441'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.4180'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").Value = '''0.9990'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '''102.82'
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").Value = '''7.129'
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").Value = '''9.177'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '''918.37'
$ws.Range("E49").Value = '  -4.49%  '
$ws.Range("D50").Value = '''34.60'
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").Value = '''0.05689'
$ws.Range("E51").Value = '  +0.77%  '
